# Commit: "Try some models for e-phys data."
# Adds 14 new rows (23-36) of model-comparison results to the "all" sheet,
# widens column B to fit the new longer kernel names, and leaves the
# selection on the last cell touched (H36), matching the authored edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("all")

# --- widen column B (new kernel names like "log_peak_green" are longer) ---
$ws.Columns.Item(2).ColumnWidth = 19.142857142857142

# --- new data rows -------------------------------------------------------
$rows = @(
    @{ r=23; A="structure";   B="peak_green_ratio"; C=-52.216850999999998; D=-47.326169;         E=0.81840000000000002; F=0.61080000000000001; H="Hypers(var_n=0.28818859003558239, var_p=0.5092397913875174)" },
    @{ r=24; A="52Structure"; B="peak_green_ratio"; C=-53.476455000000001; D=-47.712316999999999; E=0.81920000000000004; F=0.61970000000000003; H="Hypers(var_n=0.28306888121301238, ell=40.98076920387328)" },
    @{ r=25; A="32Structure"; B="peak_green_ratio"; C=-53.909047000000001; D=-48.040303000000002; E=0.81769999999999998; F=0.62560000000000004; H="Hypers(var_n=0.2820590699356596, ell=47.401370272211324)" },
    @{ r=26; A="SEStructure"; B="peak_green_ratio"; C=-53.100591000000001; D=-47.558754999999998; E=0.81989999999999996; F=0.61519999999999997; H="Hypers(var_n=0.2832017968767227, sigma_f=0.93106420538606827, ell=30.554936687641732)" },
    @{ r=27; A="structure";   B="ss_green";         C=-52.030925000000003; D=-46.109687000000001; E=0.82630000000000003; F=0.59260000000000002; H="Hypers(var_n=0.26267914889336491, var_p=0.71056386756239298)" },
    @{ r=28; A="structure";   B="ss_red";           C=-60.562184000000002; D=-55.816763000000002; E=0.73670000000000002; F=0.53039999999999998; H="Hypers(var_n=0.38732104388240945, var_p=0.73591258993921171)" },
    @{ r=29; A="structure";   B="ss_cyan";          C=-71.452267000000006; D=-69.639658999999995; E=0.44;                 F=0.26369999999999999; H="Hypers(var_n=0.65951543520724643, var_p=0.54625824498478415)" },
    @{ r=30; A="structure";   B="ss_teal";          C=-76.117149999999995; D=-76.108442999999994; E=0.0184;               F=-0.0231;              H="Hypers(var_n=0.97511162004180296, var_p=0.0067092063747461808)" },
    @{ r=31; A="structure";   B="ss_violet";        C=-73.211847000000006; D=-70.113450999999998; E=0.4476;               F=0.29139999999999999; H="Hypers(var_n=0.66322939132599756, var_p=0.84077003077915213)" },
    @{ r=32; A="structure";   B="ss_blue";          C=-70.086590999999999; D=-67.794746000000004; E=0.51239999999999997; F=0.29899999999999999; H="Hypers(var_n=0.62857027885263994, var_p=0.51022465441724618)" },
    @{ r=33; A="32Structure"; B="ss_blue";          C=-69.746010999999996; D=-67.249917999999994; E=0.5141;               F=0.34379999999999999; H="Hypers(var_n=0.38433358901934611, ell=18.260962441293746)" },
    @{ r=34; A="SEStructure"; B="ss_blue";          C=-69.248998;           D=-67.018426000000005; E=0.52100000000000002; F=0.35639999999999999; H="Hypers(var_n=0.45573304885661042, sigma_f=0.76291431833797352, ell=12.39707193305734)" },
    @{ r=35; A="32Structure"; B="log_peak_max";     C=-72.331029000000001; D=-68.688830999999993; E=0.41339999999999999; F=0.30399999999999999; H="Hypers(var_n=0.076694701182353414, ell=7.5347836014633041)" },
    @{ r=36; A="structure";   B="log_peak_green";   C=-65.045062000000001; D=-60.974125999999998; E=0.67069999999999996; F=0.52200000000000002; H="Hypers(var_n=0.48766261352288987, var_p=0.63838796083791727)" }
)

# Fill order matters for how new entries land in the shared-string table:
# rows 23-28 and 35-36 were typed row-by-row (B then H on the same row), but
# rows 29-34 had the whole B column filled first (ss_cyan/teal/violet/blue,
# blue, blue) and then the whole H column filled right after - reproduce
# that three-phase order exactly.
foreach ($row in $rows) {
    $r = $row.r
    $ws.Range("A$r").Value = $row.A
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value = $row.D
    $ws.Range("E$r").Value = $row.E
    $ws.Range("F$r").Value = $row.F

    if ($r -lt 29) {
        $ws.Range("B$r").Value = $row.B
        $ws.Range("H$r").Value = $row.H
    }
}
foreach ($row in $rows) {
    $r = $row.r
    if ($r -ge 29 -and $r -le 34) {
        $ws.Range("B$r").Value = $row.B
    }
}
foreach ($row in $rows) {
    $r = $row.r
    if ($r -ge 29 -and $r -le 34) {
        $ws.Range("H$r").Value = $row.H
    }
}
foreach ($row in $rows) {
    $r = $row.r
    if ($r -gt 34) {
        $ws.Range("B$r").Value = $row.B
        $ws.Range("H$r").Value = $row.H
    }
}

# --- leave selection where the author's cursor ended up ------------------
$ws.Range("H36").Select()
